$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.587.86"
$ws.Range("E2").Value = "  +0.67%  "

$ws.Range("D3").Value = "1.628.92"
$ws.Range("E3").Value = "  +0.99%  "

$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.62%  "

$ws.Range("D5").Value = "'213.62"
$ws.Range("E5").Value = "  -0.13%  "

$ws.Range("E6").Value = "  -0.44%  "

$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.71%  "

$ws.Range("E8").Value = "  +0.69%  "

$ws.Range("E9").Value = "  +0.13%  "

$ws.Range("D10").Value = "'19.28"
$ws.Range("E10").Value = "  +0.75%  "

$ws.Range("D11").Value = "'0.0856"
$ws.Range("E11").Value = "  -0.13%  "

$ws.Range("D12").Value = "1.856.59"
$ws.Range("E12").Value = "  +0.86%  "

$ws.Range("D13").Value = "1.630.13"
$ws.Range("E13").Value = "  -0.48%  "

$ws.Range("D14").Value = "'4.06"
$ws.Range("E14").Value = "  +0.81%  "

$ws.Range("E15").Value = "  +1.15%  "

$ws.Range("D16").Value = "'63.97"
$ws.Range("E16").Value = "  -0.24%  "

$ws.Range("D17").Value = "26.583.58"
$ws.Range("E17").Value = "  +0.58%  "

$ws.Range("D18").Value = "'234.03"
$ws.Range("E18").Value = "  +2.03%  "

$ws.Range("E19").Value = "  +2.41%  "

$ws.Range("E20").Value = "  -0.18%  "

$ws.Range("E21").Value = "  -0.40%  "

$ws.Range("E22").Value = "  -0.17%  "

$ws.Range("E23").Value = "  +1.81%  "

$ws.Range("E24").Value = "  +0.90%  "

$ws.Range("D25").Value = "'146.21"
$ws.Range("E25").Value = "  +0.37%  "

$ws.Range("E26").Value = "  -0.56%  "

$ws.Range("E27").Value = "  +1.28%  "

$ws.Range("E28").Value = "  +0.05%  "

$ws.Range("D29").Value = "'15.70"
$ws.Range("E29").Value = "  +1.27%  "

$ws.Range("D30").Value = "'0.0496"
$ws.Range("E30").Value = "  +0.10%  "

$ws.Range("E31").Value = "  -0.83%  "

$ws.Range("D32").Value = "1.526.11"
$ws.Range("E32").Value = "  +4.94%  "

$ws.Range("E33").Value = "  +0.89%  "

$ws.Range("E34").Value = "  +1.37%  "

$ws.Range("E35").Value = "  +3.80%  "

$ws.Range("E36").Value = "  -0.64%  "

$ws.Range("D37").Value = "'0.571"
$ws.Range("E37").Value = "  +0.96%  "

$ws.Range("E38").Value = "  +0.24%  "

$ws.Range("D39").Value = "'0.839"
$ws.Range("E39").Value = "  +0.78%  "

$ws.Range("E40").Value = "  -0.03%  "

$ws.Range("D41").Value = "'0.999"
$ws.Range("E41").Value = "  -0.54%  "

$ws.Range("E42").Value = "  +1.15%  "

$ws.Range("D43").Value = "1.768.44"
$ws.Range("E43").Value = "  +0.71%  "

$ws.Range("D44").Value = "'63.35"
$ws.Range("E44").Value = "  +3.24%  "

$ws.Range("D45").Value = "'0.763"
$ws.Range("E45").Value = "  +0.06%  "

$ws.Range("E46").Value = "  -3.62%  "

$ws.Range("D47").Value = "'90.09"
$ws.Range("E47").Value = "  +2.26%  "

$ws.Range("E48").Value = "  +2.02%  "

$ws.Range("E49").Value = "  +1.54%  "

$ws.Range("D50").Value = "'0.0501"
$ws.Range("E50").Value = "  -0.08%  "

$ws.Range("D51").Value = "'0.0967"
$ws.Range("E51").Value = "  +1.21%  "
